# Scheduled runner update: refresh market-board derived values (average
# price / profit columns H-N) across the per-job leve profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 2600.5
$ws.Range("J6").Value = 2600.5
$ws.Range("L6").Value = 7801.5
$ws.Range("N6").Value = -8025.5

$ws.Range("H8").Value = 469.8
$ws.Range("J8").Value = 2000
$ws.Range("L8").Value = 6000
$ws.Range("N8").Value = -6278

$ws.Range("H34").Value = 619.8
$ws.Range("I34").Value = 619.8
$ws.Range("K34").Value = 619.8
$ws.Range("M34").Value = -416.8

$ws.Range("H36").Value = 619.8
$ws.Range("I36").Value = 619.8
$ws.Range("K36").Value = 619.8
$ws.Range("M36").Value = 95.20000000000005

$ws.Range("H38").Value = 831.2
$ws.Range("I38").Value = 831.2
$ws.Range("K38").Value = 2493.6
$ws.Range("M38").Value = -2121.6

$ws.Range("H74").Value = 5528.75

$ws.Range("H77").Value = 5528.75

$ws.Range("H125").Value = 1452.5
$ws.Range("J125").Value = 1388.5714
$ws.Range("L125").Value = 12497.1426
$ws.Range("N125").Value = -17417.1426

$ws.Range("H137").Value = 1397.72
$ws.Range("I137").Value = 1104.2632
$ws.Range("K137").Value = 3312.7896
$ws.Range("M137").Value = -762.7896000000001

$ws.Range("H138").Value = 1821.1562
$ws.Range("I138").Value = 1614
$ws.Range("K138").Value = 4842
$ws.Range("M138").Value = 298

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4001.3076
$ws.Range("I32").Value = 2151.795
$ws.Range("J32").Value = 9549.846
$ws.Range("K32").Value = 2151.795
$ws.Range("L32").Value = 9549.846
$ws.Range("M32").Value = -1864.795
$ws.Range("N32").Value = -10123.846

$ws.Range("H61").Value = 3109.6155
$ws.Range("I61").Value = 2204.818
$ws.Range("K61").Value = 2204.818
$ws.Range("M61").Value = -1992.818

$ws.Range("H132").Value = 1283.4849
$ws.Range("I132").Value = 1015.23334
$ws.Range("K132").Value = 3045.70002
$ws.Range("M132").Value = -515.7000200000002

$ws.Range("H136").Value = 3109.6155
$ws.Range("I136").Value = 2204.818
$ws.Range("K136").Value = 6614.454000000001
$ws.Range("M136").Value = -4064.454000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 1000
$ws.Range("I5").Value = 1000
$ws.Range("K5").Value = 1000
$ws.Range("M5").Value = -887

$ws.Range("H20").Value = 1938.4117
$ws.Range("I20").Value = 1889.5333
$ws.Range("K20").Value = 1889.5333
$ws.Range("M20").Value = -1642.5333

$ws.Range("H134").Value = 9065.596
$ws.Range("I134").Value = 9216.147000000001
$ws.Range("K134").Value = 27648.441
$ws.Range("M134").Value = -25113.441

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 708.55554
$ws.Range("I22").Value = 255.6
$ws.Range("K22").Value = 255.6
$ws.Range("M22").Value = 94.40000000000001

$ws.Range("H94").Value = 932.8
$ws.Range("I94").Value = 755
$ws.Range("K94").Value = 755
$ws.Range("M94").Value = -304

$ws.Range("H107").Value = 385.0476
$ws.Range("I107").Value = 336.73685
$ws.Range("K107").Value = 336.73685
$ws.Range("M107").Value = 1583.26315

$ws.Range("H134").Value = 1612.6364
$ws.Range("J134").Value = 3866.3333
$ws.Range("L134").Value = 11598.9999
$ws.Range("N134").Value = -16668.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 11465.909
$ws.Range("I4").Value = 11465.909
$ws.Range("K4").Value = 34397.727
$ws.Range("M4").Value = -34285.727

$ws.Range("H6").Value = 80.2
$ws.Range("I6").Value = 50.25
$ws.Range("K6").Value = 150.75
$ws.Range("M6").Value = -37.75

$ws.Range("H122").Value = 720.61536
$ws.Range("J122").Value = 801.5
$ws.Range("L122").Value = 7213.5
$ws.Range("N122").Value = -12113.5

$ws.Range("H129").Value = 61410.668
$ws.Range("I129").Value = 681.75
$ws.Range("J129").Value = 91775.125
$ws.Range("K129").Value = 2045.25
$ws.Range("L129").Value = 275325.375
$ws.Range("M129").Value = 2954.75
$ws.Range("N129").Value = -285325.375

$ws.Range("H131").Value = 13531.072
$ws.Range("J131").Value = 14836.8
$ws.Range("L131").Value = 44510.39999999999
$ws.Range("N131").Value = -54590.39999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 5676000.5
$ws.Range("J12").Value = 2380002.8
$ws.Range("L12").Value = 2380002.8
$ws.Range("N12").Value = -2380282.8

$ws.Range("H31").Value = 1269
$ws.Range("I31").Value = 1269
$ws.Range("K31").Value = 1269
$ws.Range("M31").Value = -977

$ws.Range("H37").Value = 1269
$ws.Range("I37").Value = 1269
$ws.Range("K37").Value = 1269
$ws.Range("M37").Value = -992

$ws.Range("H80").Value = 2465
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 2465
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 2465
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -4461

$ws.Range("H83").Value = 2465
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 2465
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 12325
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -22309

$ws.Range("H102").Value = 2418.75
$ws.Range("J102").Value = 2257.75
$ws.Range("L102").Value = 2257.75
$ws.Range("N102").Value = -5501.75

$ws.Range("H122").Value = 1785.7142
$ws.Range("I122").Value = 1125
$ws.Range("J122").Value = 2666.6667
$ws.Range("K122").Value = 3375
$ws.Range("L122").Value = 8000.000100000001
$ws.Range("M122").Value = -925
$ws.Range("N122").Value = -12900.0001

$ws.Range("H126").Value = 5144282
$ws.Range("I126").Value = 11114334
$ws.Range("K126").Value = 33343002
$ws.Range("M126").Value = -33340532

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 4872.25
$ws.Range("I68").Value = 4795.8
$ws.Range("K68").Value = 4795.8
$ws.Range("M68").Value = -4046.8

$ws.Range("H71").Value = 4872.25
$ws.Range("I71").Value = 4795.8
$ws.Range("K71").Value = 23979
$ws.Range("M71").Value = -20235

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()

$ws.Range("H19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").ClearContents()

$ws.Range("H132").Value = 6216.2573
$ws.Range("I132").Value = 1214.5
$ws.Range("J132").Value = 8825.869000000001
$ws.Range("K132").Value = 3643.5
$ws.Range("L132").Value = 26477.607
$ws.Range("M132").Value = -1113.5
$ws.Range("N132").Value = -31537.607

$ws.Range("H136").Value = 23150896
$ws.Range("I136").Value = 30866838
$ws.Range("K136").Value = 92600514
$ws.Range("M136").Value = -92597964
